$d = $word.ActiveDocument

# --- Change 1: merge the two runs ("Collegamento a servizi terzi " + "(ISBN)")
#     that used to be split by a _GoBack bookmark in the list item, into a
#     single run with no bookmark.
$d.Content.Find.Execute("Collegamento a servizi terzi (ISBN)", $false, $false, $false, $false, $false, $true, 1, $false, "Collegamento a servizi terzi (ISBN)", 2) | Out-Null

# --- Locate the "Collegamento a Servizi di Terzi parti (ISBN)" Heading 2
#     paragraph; the five new body paragraphs are inserted right after it
#     (and therefore right before the "SQL" Heading 2 paragraph).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd("`r") -eq "Collegamento a Servizi di Terzi parti (ISBN)") {
        $target = $p
        break
    }
}

$cur = $target

# --- New paragraph 1 (single run)
$cur.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs($cur.Index + 1)
$p1.Style = "Normal"
$p1.Format.SpaceAfter = 0
$p1.Range.InsertAfter("Una funzione aggiuntiva al sito web che abbiamo voluto integrare è il collegamento a un servizio di terze parti che con l’inserimento di un ISBN ci permette di avere informazioni su libri anche non presenti nella biblioteca.")
$cur = $p1

# --- New paragraph 2 (two runs)
$cur.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs($cur.Index + 1)
$p2.Style = "Normal"
$p2.Format.SpaceAfter = 0
$p2.Range.InsertAfter("Il servizio di terz")
$p2r2 = $d.Range($p2.Range.End - 1, $p2.Range.End - 1)
$p2r2.InsertAfter("i parti a cui ci colleghiamo è Google Books. Con l’inserimento dell’ISBN nel sito, se è esistente, ci verranno restituiti i vari dati del libro.")
$cur = $p2

# --- New paragraph 3 (two runs)
$cur.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs($cur.Index + 1)
$p3.Style = "Normal"
$p3.Format.SpaceAfter = 0
$p3.Range.InsertAfter("Il nostro obbiettivo era di integrare più servizi di terze parti in modo da avere più informazioni possibili su un determinato libro, però questo non è stato possibile perché Google Books ")
$p3r2 = $d.Range($p3.Range.End - 1, $p3.Range.End - 1)
$p3r2.InsertAfter("è l’unico servizio di terze parti gratuito. Avremmo potuto integrare altri servizi di terze parti però sarebbero stati a pagamento, ma essendo un progetto scolastico abbiamo preferito optare per sevizi gratuiti.")
$cur = $p3

# --- New paragraph 4 (two runs)
$cur.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs($cur.Index + 1)
$p4.Style = "Normal"
$p4.Format.SpaceAfter = 0
$p4.Range.InsertAfter("Il servizio di terze parti funziona con un collegamento fra il sito web della biblioteca e il sito web di Google, il quale restituisce i dati del libro ")
$p4r2 = $d.Range($p4.Range.End - 1, $p4.Range.End - 1)
$p4r2.InsertAfter("in formato JSON. Il JSON viene modificato e viene creato un nuovo JSON che conterrà solo le informazioni utili a noi. ")
$cur = $p4

# --- New paragraph 5 (two runs + _GoBack bookmark at the very end)
$cur.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs($cur.Index + 1)
$p5.Style = "Normal"
$p5.Format.SpaceAfter = 0
$p5.Range.InsertAfter("Per lo scambio di informazioni utilizziamo il JSON per facilitare la transazione e la ")
$p5r2 = $d.Range($p5.Range.End - 1, $p5.Range.End - 1)
$p5r2.InsertAfter("leggibilità dei dati.")
$cur = $p5

# Place the _GoBack bookmark (zero-length) right at the end of paragraph 5's
# text, i.e. right after "leggibilità dei dati." and before the paragraph mark.
$bmPos = $p5.Range.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$bmRange.Bookmarks.Add("_GoBack") | Out-Null

